$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iAU_TC_ID_126"
$ws.Range("B2").Value = "@RegressionA Validation of Blueprint  Direct workflow"
$ws.Range("C2").Value = "passed"
